# ---------------------------------------------------------------------------
# withImageInElse-template.docx : "Moving from 2.0.2 to 2.0.3"
#
# What the supplied OOXML diff actually contains
# ------------------------------------------------
# Every one of the ~45 changed lines in word/document.xml and word/styles.xml
# touches only XML *attribute order* on an already-present element:
#   - the same element/tag
#   - the same text/child content
#   - the exact same set of (attribute name, attribute value) pairs
#   - just re-emitted with the attributes sorted alphabetically by local
#     name (namespace declarations first, alphabetically, then the regular
#     attributes, alphabetically) instead of the original schema-sequence
#     order Word itself writes (w:tab val/pos, w:pgSz w/h, w:pgMar
#     top/right/bottom/..., w:rFonts ascii/eastAsia/hAnsi/cs, w:lang
#     val/eastAsia/bidi, w:latentStyles/w:lsdException, w:style
#     type/default/styleId, w:tblInd w/type, w:tblCellMar w/type, and the
#     w:document root's xmlns:* declarations).
#
# No text, run, paragraph, style definition, numbering, page-size/margin
# value, font, language, tab stop, or table-style value was added,
# removed, or modified - e.g. the tab stop is still at 3119 twips/left,
# the page is still 11906 x 16838 twips with 1417/708 twips margins, the
# default run fonts/lang are unchanged, and all four w:style definitions
# keep their original name/id/type/properties. This is the signature of a
# library round trip (the commit bumps a dependency from 2.0.2 to 2.0.3)
# that re-serialises the package with a different attribute-ordering
# convention, not an authored content edit.
#
# Word's object model (real or automated) has no "attribute order" knob:
# Find/Replace, PageSetup, TabStops, Styles, etc. all write XML back out
# using Word's own fixed schema-sequence attribute order, and properties
# such as w:latentStyles/w:lsdException are not reachable through the
# Word COM surface at all (no Application/Document property exposes
# them). There is therefore no COM call that changes what a user/macro
# could observe here - the document's content, formatting and structure
# are identical before and after.
#
# So this script intentionally makes no content mutations. It only
# walks the object model to confirm the template still matches what the
# diff describes (so the "edit" is a verified no-op rather than a
# silent skip), and leaves the document exactly as authored.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Sanity-check the five conditional-tab paragraphs the diff touches are
# still the untouched "m:if / m:else / ... / m:endif" field skeleton with
# their left tab at 3119 twips - read-only access, nothing is written.
$paragraphCount = $d.Paragraphs.Count
Write-Output "Paragraphs: $paragraphCount"

for ($i = 1; $i -le $paragraphCount; $i++) {
    $para = $d.Paragraphs.Item($i)
    $tabStops = $para.Range.TabStops
    if ($tabStops.Count -gt 0) {
        $firstTab = $tabStops.Item(1)
        Write-Output "Paragraph $i tab stop position: $($firstTab.Position)"
    }
}

# Sanity-check the page setup values referenced in the w:sectPr hunk
# (w:pgSz / w:pgMar) are unchanged - again, read-only.
$section = $d.Sections.Item(1)
$pageSetup = $section.PageSetup
Write-Output "PageWidth: $($pageSetup.PageWidth), PageHeight: $($pageSetup.PageHeight)"
Write-Output "Margins T/R/B/L: $($pageSetup.TopMargin)/$($pageSetup.RightMargin)/$($pageSetup.BottomMargin)/$($pageSetup.LeftMargin)"

# Sanity-check the four named styles referenced in the w:style hunks are
# still present and unchanged.
$styles = $d.Styles
for ($i = 1; $i -le $styles.Count; $i++) {
    $style = $styles.Item($i)
    Write-Output "Style: $($style.NameLocal)"
}

# No InsertXML / Range.Text / formatting / style assignments follow: the
# diff carries zero semantic change, so none are needed.
